$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: the "Price" (D) column holds numeric-looking text (e.g. "244.45",
# "0.9450" with a significant trailing zero). Excel's Range.Value setter
# auto-converts numeric-looking strings to real numbers, which would both
# change the cell's stored type and silently drop meaningful trailing
# zeros (e.g. "6.410" -> 6.41). Forcing the cell to Text format before
# assigning keeps the literal string intact, matching the source data.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Rows 2-8: price-only updates
Set-TextValue $ws.Range("D2") "244.45"
Set-TextValue $ws.Range("D3") "21.83"
Set-TextValue $ws.Range("D4") "5.388"
Set-TextValue $ws.Range("D5") "0.06001"
Set-TextValue $ws.Range("D6") "3.374"
Set-TextValue $ws.Range("D7") "0.8139"
Set-TextValue $ws.Range("D8") "0.9450"

# Row 9: One -> WazirX
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D9") "0.1428"
$ws.Range("E9").Value = "8WazirXWRX"

# Row 10: WazirX -> MandalaExchangeToken
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D10") "0.07397"
$ws.Range("E10").Value = "9MandalaExchangeTokenMDX"

# Row 11: MandalaExchangeToken -> LiechtensteinCryptoassetsExchange
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D11") "0.03333"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"

# Row 12: LiechtensteinCryptoassetsExchange -> BitrueCoin
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D12") "0.03052"
$ws.Range("E12").Value = "11BitrueCoinBTR"

# Row 13: BitrueCoin -> BitMartToken
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D13") "0.09405"
$ws.Range("E13").Value = "12BitMartTokenBMX"

# Row 14: BitMartToken -> MCDex
$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws.Range("D14") "4.003"
$ws.Range("E14").Value = "13MCDexMCB"

# Row 15: MCDex -> BitForexToken
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D15") "0.001591"
$ws.Range("E15").Value = "14BitForexTokenBF"

# Row 16: BitForexToken -> CoinExToken
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws.Range("D16") "0.04817"
$ws.Range("E16").Value = "15CoinExTokenCET"

# Row 17: CoinExToken -> One
$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D17") "0.0005912"
$ws.Range("E17").Value = "16OneONEWorstin24h"

# Price-only updates further down
Set-TextValue $ws.Range("D18") "0.006249"
Set-TextValue $ws.Range("D19") "0.005001"
Set-TextValue $ws.Range("D20") "0.0009866"
Set-TextValue $ws.Range("D22") "3.681"
Set-TextValue $ws.Range("D23") "6.410"
Set-TextValue $ws.Range("D24") "2.219"

# Row 27: price + volume label
Set-TextValue $ws.Range("D27") "0.0002451"
$ws.Range("E27").Value = "26UpBotsUBXTBestin24h"

# Rows 41-45: price-only updates
Set-TextValue $ws.Range("D41") "0.006492"
Set-TextValue $ws.Range("D42") "0.1073"
Set-TextValue $ws.Range("D43") "0.003101"
Set-TextValue $ws.Range("D44") "0.005244"
Set-TextValue $ws.Range("D45") "0.00005265"

# Row 47: price + volume label
Set-TextValue $ws.Range("D47") "0.9692"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

# Row 48: price + volume label
Set-TextValue $ws.Range("D48") "0.01516"
$ws.Range("E48").Value = "47BOLOBOLO"
